# Add the 2022-Q3 quarterly sheet, insert it into the "总计" (total) summary
# sheet, and re-order worksheets so the new quarter sits right after the
# summary sheet and before the older quarters.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" summary sheet: insert a new row for 2022-Q3 right under the
#    header row, pushing the existing 2021-Q4 / 2021-Q3 rows down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.01

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2

# ---------------------------------------------------------------------
# 2. New "2022-Q3" detail sheet: clone the "2021-Q4" sheet (same layout /
#    styling) right after "总计", rename it, and replace its data.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($null, $total)

$new = $wb.Worksheets.Item(2)
$new.Name = "2022-Q3"

$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"

$data = @(
  @(0, "000804", "中信建投稳利混合A", "0.21", "38.27", "1.70", "0.0036", 10),
  @(1, "003308", "中信建投睿利灵活配置混合A", "0.07", "93.78", "3.85", "0.0027", 8),
  @(2, "006844", "中信建投稳利混合C", "0.10", "38.27", "1.70", "0.0017", 10),
  @(3, "004635", "中信建投睿利灵活配置混合C", "0.03", "93.78", "3.85", "0.0012", 8)
)

$r = 2
foreach ($row in $data) {
    $new.Cells.Item($r, 1).Value = $row[0]
    $new.Cells.Item($r, 2).Value = "'" + $row[1]
    $new.Cells.Item($r, 3).Value = $row[2]
    $new.Cells.Item($r, 4).Value = "'" + $row[3]
    $new.Cells.Item($r, 5).Value = "'" + $row[4]
    $new.Cells.Item($r, 6).Value = "'" + $row[5]
    $new.Cells.Item($r, 7).Value = "'" + $row[6]
    $new.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# The cloned sheet had 5 data rows; 2022-Q3 only has 4, drop the leftover.
$new.Rows.Item(6).Delete()
